$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (columns H..AA) per updated election results
$ws.Range("H2").Value  = 49
$ws.Range("I2").Value  = 136
$ws.Range("J2").Value  = 525
$ws.Range("K2").Value  = 2
$ws.Range("L2").Value  = 121
$ws.Range("M2").Value  = 12
$ws.Range("N2").Value  = 89
$ws.Range("O2").Value  = 1
$ws.Range("P2").Value  = 1
$ws.Range("Q2").Value  = 1
$ws.Range("R2").Value  = 6
$ws.Range("S2").Value  = 58
$ws.Range("T2").Value  = 90
$ws.Range("U2").Value  = 7
$ws.Range("V2").Value  = 844
$ws.Range("W2").Value  = 0
$ws.Range("X2").Value  = 774
$ws.Range("Y2").Value  = 0
$ws.Range("Z2").Value  = 15
$ws.Range("AA2").Value = 10
